$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
    @(1.02, 1.078796209455121, 1.065371551027789, 1.079562070010494, 1.084260198779201, 1.042938977784816, 1.083683902046142, 1.068084970606566, 1.082237796409068, 1.086923680630766, 1.085222857487328),
    @(1.02, 1.081365451188037, 1.066607333607261, 1.081721707774545, 1.086176982265449, 1.04327750284408, 1.085906404367815, 1.069135088482449, 1.08421227103245, 1.088656794800199, 1.087448516017222),
    @(1.02, 1.083020839256943, 1.067402283777456, 1.083112673978134, 1.087411109476947, 1.043493374184057, 1.087337427875716, 1.069809345860253, 1.08548305714028, 1.089771654960703, 1.088881571742646),
    @(1.02, 1.083715115140304, 1.067735377492369, 1.083695926238086, 1.087928493547466, 1.043583372508979, 1.087937374468997, 1.070091566713916, 1.086015695845035, 1.090238796876561, 1.089482370328904),
    @(1.02, 1.083831591742855, 1.067791241210802, 1.083793769398386, 1.088015280904546, 1.043598439599718, 1.088038012212919, 1.070138880800774, 1.086105035442126, 1.090317142311444, 1.089583150989965),
    @(1.02, 1.083030122623678, 1.067406738903143, 1.083120473306175, 1.087418028415809, 1.043494579699833, 1.087345450843356, 1.069813121746203, 1.085490180523141, 1.089777902965455, 1.08888960610382),
    @(1.02, 1.079665991512102, 1.065790170254847, 1.080293290844371, 1.084909279391052, 1.043054045739554, 1.084436498342091, 1.06844095815391, 1.082906515338492, 1.087510775502377, 1.085976522556342),
    @(1.02, 1.073681604954219, 1.062904879800143, 1.075260224723335, 1.080439926693055, 1.042253136059071, 1.079254468153236, 1.065982117983996, 1.078299833772769, 1.083464039480723, 1.080787133290232),
    @(1.02, 1.069651141306174, 1.06095552566747, 1.071868008154136, 1.077425635949529, 1.041702208226538, 1.075759489038284, 1.064314268166859, 1.07519016616819, 1.080729474193555, 1.077287190904057),
    @(1.02, 1.067895539065916, 1.060105045467958, 1.070389859488109, 1.076111719392015, 1.041459524244771, 1.074235984393459, 1.063585029116225, 1.073833987240298, 1.079536225403762, 1.075761522707888),
    @(1.02, 1.067241811524915, 1.059788155241618, 1.069839365253017, 1.075622323842786, 1.041368751840932, 1.073668510411051, 1.063313076137782, 1.07332874394442, 1.079091584786385, 1.075193242847334),
    @(1.02, 1.067382112504258, 1.059856174173485, 1.069957514250337, 1.075727362521552, 1.041388251421268, 1.073790307674652, 1.063371460258067, 1.073437188868386, 1.079187026355784, 1.07531521307702),
    @(1.02, 1.067841535107213, 1.060078871431032, 1.070344385197434, 1.076071293545634, 1.041452033846032, 1.07418910927509, 1.063562571582298, 1.073792254494907, 1.079499500358221, 1.075714581021474),
    @(1.02, 1.068124384272536, 1.060215951475875, 1.070582556380112, 1.076283020908102, 1.041491248728073, 1.07443461366294, 1.063680177648589, 1.07401082228033, 1.079691837264787, 1.075960434053714),
    @(1.02, 1.069767431081724, 1.061011832513132, 1.071965908129554, 1.077512649380804, 1.041718226657635, 1.075860380747646, 1.064362514970747, 1.075279963700342, 1.080808469876642, 1.077388225891218),
    @(1.02, 1.07079524968485, 1.061509338960968, 1.072831126092479, 1.078281603233, 1.041859492635393, 1.076751972310459, 1.064788624463698, 1.076073441609608, 1.08150642389837, 1.078281083616317),
    @(1.02, 1.071393760904259, 1.061798910410987, 1.073334899356113, 1.078729283479559, 1.041941492923641, 1.077271046192388, 1.065036488303172, 1.076535333686661, 1.081912647092098, 1.078800894642697),
    @(1.02, 1.071597670240845, 1.061897543055188, 1.073506522845039, 1.078881790015458, 1.041969385688961, 1.077447872879995, 1.065120889022843, 1.076692670522273, 1.082051010314012, 1.078977972444483),
    @(1.02, 1.07068507812163, 1.061456025071037, 1.072738389176228, 1.078199188678589, 1.04184437733672, 1.076656414330758, 1.064742977298954, 1.075988405472673, 1.081431631496501, 1.078185389933324),
    @(1.02, 1.067706291760919, 1.060013319989344, 1.070230501592181, 1.075970051980253, 1.041433268937818, 1.07407171607593, 1.06350632408456, 1.073687738236354, 1.079407523920267, 1.075597021110505),
    @(1.02, 1.065824015429978, 1.0591005291463, 1.068645314443235, 1.074560683207125, 1.041171146578743, 1.072437466546985, 1.062722524798447, 1.072232527242017, 1.078126681621149, 1.073960450759779),
    @(1.02, 1.06682275611039, 1.059584965383454, 1.069486462389386, 1.075308571279914, 1.04131045072131, 1.073304697399482, 1.063138633097936, 1.073004801079237, 1.078806471070562, 1.074828913179558),
    @(1.02, 1.070734862989811, 1.0614801172245, 1.072780295756381, 1.078236430854462, 1.041851208518322, 1.076699595883257, 1.064763605383886, 1.076026832514781, 1.081465429685914, 1.078228632808581),
    @(1.02, 1.075235694460167, 1.063655262380905, 1.076567706755553, 1.081601329960672, 1.04246315285916, 1.080601056527492, 1.066622748955294, 1.079497395341024, 1.084516549485469, 1.082135633974434)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $colLetter = $cols[$j]
        $ws.Range("$colLetter$rowNum").Value = $rowVals[$j]
    }
}